$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "phone" looks numeric (leading zero) so force it to stay text before
# writing the value, otherwise Excel's auto-detection would coerce it to a
# number and drop the leading zero.
$ws.Range("D3").NumberFormat = "@"

$ws.Range("A3").Value = "SL-20251128-001"
$ws.Range("B3").Value = "2025-11-28 00:21:24"
$ws.Range("C3").Value = "Fahad Ahmed"
$ws.Range("D3").Value = "0502992932"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 175
$ws.Range("G3").Value = 175
$ws.Range("H3").Value = "pending"
# payment_intent_id / redirect_url are blank text cells (present but empty),
# not simply missing -- a leading quote forces Excel to keep them as an
# actual (empty) text value instead of clearing the cell outright.
$ws.Range("I3").Value = "'"
$ws.Range("J3").Value = "error"
$ws.Range("K3").Value = "'"
$ws.Range("L3").Value = 50
